# "unfertig!!! Live Pie Chart" - fill in the (until now empty) "D" column
# of the Gewichte sheet with weight=1 for every row, and backfill a couple
# of missing "C" column values that were left blank. Also row 13's C cell
# held the stray shared string "s" - replace it with the numeric weight.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichte")

$xlPasteFormats = -4122

# Rows 3-12: just set D to 1; C/D already share the same style in these rows.
3..12 | ForEach-Object {
    $ws.Cells.Item($_, 4).Value = 1
}

# Row 13: C13 currently holds the text "s" - replace with numeric 1, then
# give D13 the same (now-numeric) style as C13 before stamping its value.
$ws.Range("C13").Value = 1
$ws.Range("C13").Copy()
$ws.Range("D13").PasteSpecial($xlPasteFormats)
$ws.Range("D13").Value = 1

# Row 14: C14 was blank -> 1; D14 copies C14's style, then gets 1.
$ws.Range("C14").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)
$ws.Range("D14").Value = 1

# Row 15: C15 was blank -> 2; D15 copies C15's style, then gets 1.
$ws.Range("C15").Value = 2
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Value = 1

# Rows 16-18: C already populated; D copies C's style, then gets 1.
16..18 | ForEach-Object {
    $cC = $ws.Cells.Item($_, 3)
    $cD = $ws.Cells.Item($_, 4)
    $cC.Copy()
    $cD.PasteSpecial($xlPasteFormats)
    $cD.Value = 1
}

# Row 19: C19 was blank -> 2; D19 copies C19's style, then gets 1.
$ws.Range("C19").Value = 2
$ws.Range("C19").Copy()
$ws.Range("D19").PasteSpecial($xlPasteFormats)
$ws.Range("D19").Value = 1

# Row 20: C20 was blank -> 3; D20 copies C20's style, then gets 1.
$ws.Range("C20").Value = 3
$ws.Range("C20").Copy()
$ws.Range("D20").PasteSpecial($xlPasteFormats)
$ws.Range("D20").Value = 1

$excel.CutCopyMode = $false
